# Update "想去人数" (attendance interest count) figures in column F
# for the "展览" and "全部类型" worksheets, as per the source diff.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 3319
$ws1.Range("F5").Value = 6930
$ws1.Range("F6").Value = 2277
$ws1.Range("F13").Value = 165
$ws1.Range("F14").Value = 410

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 3319
$ws4.Range("F6").Value = 6930
$ws4.Range("F7").Value = 2277
$ws4.Range("F14").Value = 165
$ws4.Range("F15").Value = 410
